# Fix render text thumb issue: populate heading/subheading/extraheading
# for row 2 and correct the release date / release hour values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "op-test"
$ws.Range("F2").Value = "who is there"
$ws.Range("G2").Value = 4444

$ws.Range("L2").Value = "2023-10-13T00:00:00+00:00"
$ws.Range("M2").Value = "19:45"
